$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (overwrite existing data)
$ws.Range("A2").Value = "sbkuzh"
$ws.Range("B2").Value = "task_2050-01-01_UZH_LARGE_READY"
$ws.Range("C2").Value = "2024-07-20 00:40:59"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2050-01-01"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "LARGE"
$ws.Range("H2").Value = "READY"

# Row 3 (new)
$ws.Range("A3").Value = "sbkrzs"
$ws.Range("B3").Value = "task_2034-01-01_RZS_LARGE_ERROR"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2034-01-01"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Value = "LARGE"
$ws.Range("H3").Value = "ERROR"
$ws.Range("I3").Value = "Missing file task_2034-01-01_RZS_LARGE.xlsx"

# Row 4 (new)
$ws.Range("A4").Value = "sbkrzs"
$ws.Range("B4").Value = "task_2033-01-01_RZS_LARGE_ERROR"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2033-01-01"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").Value = "LARGE"
$ws.Range("H4").Value = "ERROR"
$ws.Range("I4").Value = "Missing file task_2033-01-01_RZS_LARGE.xlsx"

# Row 5 (new)
$ws.Range("A5").Value = "sbkzbs"
$ws.Range("B5").Value = "task_2033-01-01_ZBS_LARGE_ERROR"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2033-01-01"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").Value = "LARGE"
$ws.Range("H5").Value = "ERROR"
$ws.Range("I5").Value = "A large task is already scheduled for this date"

# Row 6 (new)
$ws.Range("A6").Value = "sbkrzs"
$ws.Range("B6").Value = "task_2041-01-01_UBS_LARGE_ERROR"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2041-01-01"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").Value = "LARGE"
$ws.Range("H6").Value = "ERROR"
$ws.Range("I6").Value = "Missing file task_2041-01-01_UBS_LARGE.xlsx"

# Row 7 (new)
$ws.Range("A7").Value = "sbkrzs"
$ws.Range("B7").Value = "task_2032-01-01_RZS_LARGE_ERROR"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2032-01-01"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").Value = "LARGE"
$ws.Range("H7").Value = "ERROR"
$ws.Range("I7").Value = "Missing file task_2032-01-01_RZS_LARGE.xlsx"

# Row 8 (new)
$ws.Range("A8").Value = "sbkzbz"
$ws.Range("B8").Value = "task_2024-07-22_ZBZ_SMALL_READY"
$ws.Range("C8").Value = "2024-07-20 00:41:56"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2024-07-22"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").Value = "SMALL"
$ws.Range("H8").Value = "READY"

# Row 9 (new)
$ws.Range("A9").Value = "sbkhsg"
$ws.Range("B9").Value = "task_2024-07-20_HSG_SMALL_DONE"
$ws.Range("C9").Value = "2024-07-20 00:42:06"
$ws.Range("E9").Value = "2024-07-20 00:42:16"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2024-07-20"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").Value = "SMALL"
$ws.Range("H9").Value = "DONE"
